# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45182 (2023-09-13) to 45184 (2023-09-15).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 77
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value2 = 45184
    }
}
